{"js": "// Update the worksheet date header and the 25 \"two-digit \u00f7 one-digit\"\n// division prompts to the new day's values. Each table cell (and the\n// title paragraph) holds exactly one run, so the text is replaced in\n// place via the paragraph, which keeps the existing run formatting\n// (font / size) untouched.\n\n// Map of the paragraph's current (\"old\") text to its new text, in the\n// same document order the paragraphs occur in (title, then each\n// populated table row left-to-right, top-to-bottom).\nconst replacements = [\n  [\"2024-11-30 Saturday\", \"2024-12-01 Sunday\"],\n  [\"30\u00f72=\", \"74\u00f72=\"],\n  [\"23\u00f77=\", \"74\u00f75=\"],\n  [\"94\u00f78=\", \"46\u00f74=\"],\n  [\"19\u00f77=\", \"35\u00f74=\"],\n  [\"92\u00f73=\", \"86\u00f78=\"],\n  [\"70\u00f78=\", \"94\u00f78=\"],\n  [\"54\u00f73=\", \"71\u00f79=\"],\n  [\"55\u00f72=\", \"29\u00f76=\"],\n  [\"76\u00f73=\", \"89\u00f75=\"],\n  [\"85\u00f73=\", \"61\u00f72=\"],\n  [\"32\u00f74=\", \"53\u00f78=\"],\n  [\"23\u00f75=\", \"83\u00f78=\"],\n  [\"74\u00f76=\", \"87\u00f74=\"],\n  [\"35\u00f72=\", \"15\u00f76=\"],\n  [\"72\u00f76=\", \"67\u00f73=\"],\n  [\"28\u00f72=\", \"38\u00f78=\"],\n  [\"16\u00f73=\", \"99\u00f73=\"],\n  [\"55\u00f76=\", \"45\u00f79=\"],\n  [\"96\u00f77=\", \"26\u00f77=\"],\n  [\"90\u00f73=\", \"88\u00f76=\"],\n  [\"69\u00f76=\", \"45\u00f75=\"],\n  [\"60\u00f79=\", \"41\u00f76=\"],\n  [\"92\u00f76=\", \"50\u00f79=\"],\n  [\"84\u00f73=\", \"94\u00f74=\"],\n  [\"11\u00f76=\", \"66\u00f74=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Only paragraphs that actually carry text (the title + the 25 filled\n// table cells) take part in the substitution; every other paragraph in\n// the table is empty and is left untouched.\nconst nonEmpty = [];\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text !== \"\") {\n    nonEmpty.push(p);\n  }\n}\n\nif (nonEmpty.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" text-bearing paragraphs, found \" + nonEmpty.length\n  );\n}\n\nfor (let i = 0; i < nonEmpty.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = nonEmpty[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \": expected \" + JSON.stringify(oldText) + \" but found \" + JSON.stringify(para.text)\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 \"two-digit \u00f7 one-digit\"\n# division prompts to the new day's values.\n#\n# Each (old, new) pair below is unique and is applied with\n# Find.Execute(..., Replace:=wdReplaceOne) against a *fresh*\n# $d.Content range every time, so the match position of one\n# replacement never affects later ones (important since one of the\n# new values, \"94\u00f78=\", equals an old value elsewhere in the sheet\n# that is handled earlier in this list).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-30 Saturday\", \"2024-12-01 Sunday\"),\n    @(\"30\u00f72=\", \"74\u00f72=\"),\n    @(\"23\u00f77=\", \"74\u00f75=\"),\n    @(\"94\u00f78=\", \"46\u00f74=\"),\n    @(\"19\u00f77=\", \"35\u00f74=\"),\n    @(\"92\u00f73=\", \"86\u00f78=\"),\n    @(\"70\u00f78=\", \"94\u00f78=\"),\n    @(\"54\u00f73=\", \"71\u00f79=\"),\n    @(\"55\u00f72=\", \"29\u00f76=\"),\n    @(\"76\u00f73=\", \"89\u00f75=\"),\n    @(\"85\u00f73=\", \"61\u00f72=\"),\n    @(\"32\u00f74=\", \"53\u00f78=\"),\n    @(\"23\u00f75=\", \"83\u00f78=\"),\n    @(\"74\u00f76=\", \"87\u00f74=\"),\n    @(\"35\u00f72=\", \"15\u00f76=\"),\n    @(\"72\u00f76=\", \"67\u00f73=\"),\n    @(\"28\u00f72=\", \"38\u00f78=\"),\n    @(\"16\u00f73=\", \"99\u00f73=\"),\n    @(\"55\u00f76=\", \"45\u00f79=\"),\n    @(\"96\u00f77=\", \"26\u00f77=\"),\n    @(\"90\u00f73=\", \"88\u00f76=\"),\n    @(\"69\u00f76=\", \"45\u00f75=\"),\n    @(\"60\u00f79=\", \"41\u00f76=\"),\n    @(\"92\u00f76=\", \"50\u00f79=\"),\n    @(\"84\u00f73=\", \"94\u00f74=\"),\n    @(\"11\u00f76=\", \"66\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        Write-Output (\"WARNING: could not find \" + $oldText)\n    }\n}\n"}
